$r = Get-Command "TotallyFakeCmdletXYZ" -ErrorAction Stop
Write-Output "Result=$r"
